$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = [double]"6.532620551785419E-08"
$ws.Range("E2").Value = [double]"6.532620551785419E-08"

# Row 3
$ws.Range("D3").Value = [double]"1.008317268527404E-46"
$ws.Range("E3").Value = [double]"1.008317268527404E-46"

# Row 4
$ws.Range("D4").Value = [double]"3.371428394512979E-15"
$ws.Range("E4").Value = [double]"3.371428394512979E-15"

# Row 5
$ws.Range("C5").Value = $true
$ws.Range("D5").Value = 0.005647965465362825
$ws.Range("E5").Value = 0.005647965465362825

# Row 6
$ws.Range("D6").Value = 0.9999996253748886
$ws.Range("E6").Value = 0.9999996253748886

# Row 7
$ws.Range("F7").Value = 2.467167377471924
$ws.Range("G7").Value = 0.8333333333333334
